# repull data, push all data, mean calculation
# Update the "dSF" column (F) values for rows where the re-pulled data differs
# from the originally stored dS0 (column E) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -11
    3  = -8
    8  = -3
    14 = -6
    16 = -11
    22 = -8
    25 = -2
    34 = 4
    48 = -3
    50 = 1
    53 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
